# Apply the "Analysis Results" sheet updates:
#  - relabel several metric names to include their units
#  - swap Starting/Ending SoC (%) values (B6/B7)
#  - swap Highest/Lowest Cell Voltage values (B16/B17)
#  - fix sign of Regenerative Effectiveness and recompute several values
#  - remove the old "Maximum BMS Temperature in C" row, shifting rows 31-42
#    up by one (with recomputed values) and append a new row 43
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column A labels (row, new label per diff) ---
$ws.Range("A8").Value = "Total distance covered (km)"
$ws.Range("A9").Value = "Total energy consumption(WH/KM)"
$ws.Range("A10").Value = "Total SOC consumed(%)"
$ws.Range("A12").Value = "Peak Power(kW)"
$ws.Range("A13").Value = "Average Power(kW)"
$ws.Range("A14").Value = "Total Energy Regenerated(kWh)"
$ws.Range("A15").Value = "Regenerative Effectiveness(%)"
$ws.Range("A16").Value = "Highest Cell Voltage(V)"
$ws.Range("A17").Value = "Lowest Cell Voltage(V)"
$ws.Range("A18").Value = "Difference in Cell Voltage(V)"
$ws.Range("A19").Value = "Minimum Temperature(C)"
$ws.Range("A20").Value = "Maximum Temperature(C)"
$ws.Range("A21").Value = "Difference in Temperature(C)"
$ws.Range("A22").Value = "Maximum Fet Temperature-BMS(C)"
$ws.Range("A23").Value = "Maximum Afe Temperature-BMS(C)"
$ws.Range("A24").Value = "Maximum PCB Temperature-BMS(C)"
$ws.Range("A25").Value = "Maximum MCU Temperature(C)"
$ws.Range("A26").Value = "Maximum Motor Temperature(C)"
$ws.Range("A27").Value = "Abnormal Motor Temperature Detected(C)"
$ws.Range("A28").Value = "highest cell temp(C)"
$ws.Range("A29").Value = "lowest cell temp(C)"
$ws.Range("A30").Value = "Difference between Highest and Lowest Cell Temperature at 100% SOC(C)"
$ws.Range("A31").Value = "Battery Voltage(V)"
$ws.Range("A32").Value = "Total energy charged(kWh)"
$ws.Range("A33").Value = "Electricity consumption units(kW)"
$ws.Range("A34").Value = "Idling time percentage"
$ws.Range("A35").Value = "Time spent in 0-10 km/h"
$ws.Range("A36").Value = "Time spent in 10-20 km/h"
$ws.Range("A37").Value = "Time spent in 20-30 km/h"
$ws.Range("A38").Value = "Time spent in 30-40 km/h"
$ws.Range("A39").Value = "Time spent in 40-50 km/h"
$ws.Range("A40").Value = "Time spent in 50-60 km/h"
$ws.Range("A41").Value = "Time spent in 60-70 km/h"
$ws.Range("A42").Value = "Time spent in 70-80 km/h"

# --- Update column B numeric values ---
$ws.Range("B6").Value = 96
$ws.Range("B7").Value = 35
$ws.Range("B15").Value = 4.266675008452464
$ws.Range("B16").Value = 3.33
$ws.Range("B17").Value = 3.082
$ws.Range("B21").Value = 7
$ws.Range("B31").Value = 55
$ws.Range("B32").Value = 1.281246388888889
$ws.Range("B33").Value = [double]"9.149145878955218e-08"
$ws.Range("B34").Value = 5.959138900315371
$ws.Range("B35").Value = 12.43384066913479
$ws.Range("B36").Value = 3.587001234060057
$ws.Range("B37").Value = 7.061565885095297
$ws.Range("B38").Value = 33.23460852872618
$ws.Range("B39").Value = 22.65734265734266
$ws.Range("B40").Value = 13.45399698340875
$ws.Range("B41").Value = 0.4607157548334019
$ws.Range("B42").Value = 0

# --- Add new row 43 (extends dimension A1:B42 -> A1:B43) ---
$ws.Range("A43").Value = "Time spent in 80-90 km/h"
$ws.Range("B43").Value = 0
